$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.758258666666666
$ws.Range("H2").Value = 5.274775999999999
$ws.Range("I2").Value = 0.1132279568112417
$ws.Range("J2").Value = 0.1132279568112417
$ws.Range("M2").Value = 43.97948166666666
$ws.Range("N2").Value = 131.938445
$ws.Range("O2").Value = 0.3260725128076164
$ws.Range("P2").Value = 0.3260725128076164
$ws.Range("Q2").Value = 77.32730479592443
$ws.Range("R2").Value = 695.9457431633199
$ws.Range("S2").Value = 0.03692052439751384
$ws.Range("T2").Value = 0.03692052439751384

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.758258666666666
$ws.Range("H3").Value = 5.274775999999999
$ws.Range("I3").Value = 0.1132279568112417
$ws.Range("J3").Value = 0.1132279568112417
$ws.Range("M3").Value = 57.80064033333333
$ws.Range("O3").Value = 0.4285452970598356
$ws.Range("P3").Value = 0.4285452970598356
$ws.Range("Q3").Value = 101.6284768049662
$ws.Range("R3").Value = 914.6562912446958
$ws.Range("S3").Value = 0.0485233083871518
$ws.Range("T3").Value = 0.0485233083871518

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.758258666666666
$ws.Range("H4").Value = 5.274775999999999
$ws.Range("I4").Value = 0.1132279568112417
$ws.Range("J4").Value = 0.1132279568112417
$ws.Range("M4").Value = 20.92900166666667
$ws.Range("N4").Value = 62.787005
$ws.Range("O4").Value = 0.15517172793733
$ws.Range("P4").Value = 0.15517172793733
$ws.Range("Q4").Value = 36.79859856509778
$ws.Range("R4").Value = 331.1873870858799
$ws.Range("S4").Value = 0.01756977770921374
$ws.Range("T4").Value = 0.01756977770921374

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 1.758258666666666
$ws.Range("H5").Value = 5.274775999999999
$ws.Range("I5").Value = 0.1132279568112417
$ws.Range("J5").Value = 0.1132279568112417
$ws.Range("M5").Value = 12.167261
$ws.Range("N5").Value = 36.501783
$ws.Range("O5").Value = 0.0902104621952179
$ws.Range("P5").Value = 0.0902104621952179
$ws.Range("Q5").Value = 21.39319210284533
$ws.Range("R5").Value = 192.538728925608
$ws.Range("S5").Value = 0.01021434631736228
$ws.Range("T5").Value = 0.01021434631736228

# Row 6
$ws.Range("I6").Value = 0.1732550390834427
$ws.Range("J6").Value = 0.1732550390834427
$ws.Range("M6").Value = 43.97948166666666
$ws.Range("N6").Value = 131.938445
$ws.Range("O6").Value = 0.3260725128076164
$ws.Range("P6").Value = 0.3260725128076164
$ws.Range("Q6").Value = 118.3218843820472
$ws.Range("R6").Value = 1064.896959438425
$ws.Range("S6").Value = 0.05649370595051995
$ws.Range("T6").Value = 0.05649370595051995

# Row 7
$ws.Range("I7").Value = 0.1732550390834427
$ws.Range("J7").Value = 0.1732550390834427
$ws.Range("M7").Value = 57.80064033333333
$ws.Range("O7").Value = 0.4285452970598356
$ws.Range("P7").Value = 0.4285452970598356
$ws.Range("S7").Value = 0.07424763219112739
$ws.Range("T7").Value = 0.07424763219112739

# Row 8
$ws.Range("I8").Value = 0.1732550390834427
$ws.Range("J8").Value = 0.1732550390834427
$ws.Range("M8").Value = 20.92900166666667
$ws.Range("N8").Value = 62.787005
$ws.Range("O8").Value = 0.15517172793733
$ws.Range("P8").Value = 0.15517172793733
$ws.Range("Q8").Value = 56.3071419123139
$ws.Range("R8").Value = 506.764277210825
$ws.Range("S8").Value = 0.02688428378842745
$ws.Range("T8").Value = 0.02688428378842744

# Row 9
$ws.Range("I9").Value = 0.1732550390834427
$ws.Range("J9").Value = 0.1732550390834427
$ws.Range("M9").Value = 12.167261
$ws.Range("N9").Value = 36.501783
$ws.Range("O9").Value = 0.0902104621952179
$ws.Range("P9").Value = 0.0902104621952179
$ws.Range("Q9").Value = 32.73465704302167
$ws.Range("R9").Value = 294.611913387195
$ws.Range("S9").Value = 0.01562941715336791
$ws.Range("T9").Value = 0.01562941715336791

# Row 10
$ws.Range("G10").Value = 9.938311666666667
$ws.Range("H10").Value = 29.814935
$ws.Range("I10").Value = 0.6400052196548212
$ws.Range("J10").Value = 0.640005219654821
$ws.Range("M10").Value = 43.97948166666666
$ws.Range("N10").Value = 131.938445
$ws.Range("O10").Value = 0.3260725128076164
$ws.Range("P10").Value = 0.3260725128076164
$ws.Range("Q10").Value = 437.0817957417861
$ws.Range("R10").Value = 3933.736161676075
$ws.Range("S10").Value = 0.208688110182838
$ws.Range("T10").Value = 0.208688110182838

# Row 11
$ws.Range("G11").Value = 9.938311666666667
$ws.Range("H11").Value = 29.814935
$ws.Range("I11").Value = 0.6400052196548212
$ws.Range("J11").Value = 0.640005219654821
$ws.Range("M11").Value = 57.80064033333333
$ws.Range("O11").Value = 0.4285452970598356
$ws.Range("P11").Value = 0.4285452970598356
$ws.Range("Q11").Value = 574.4407781655706
$ws.Range("R11").Value = 5169.967003490135
$ws.Range("S11").Value = 0.2742712269768207
$ws.Range("T11").Value = 0.2742712269768207

# Row 12
$ws.Range("G12").Value = 9.938311666666667
$ws.Range("H12").Value = 29.814935
$ws.Range("I12").Value = 0.6400052196548212
$ws.Range("J12").Value = 0.640005219654821
$ws.Range("M12").Value = 20.92900166666667
$ws.Range("N12").Value = 62.787005
$ws.Range("O12").Value = 0.15517172793733
$ws.Range("P12").Value = 0.15517172793733
$ws.Range("Q12").Value = 207.9989414355195
$ws.Range("R12").Value = 1871.990472919675
$ws.Range("S12").Value = 0.09931071582274903
$ws.Range("T12").Value = 0.099310715822749

# Row 13
$ws.Range("G13").Value = 9.938311666666667
$ws.Range("H13").Value = 29.814935
$ws.Range("I13").Value = 0.6400052196548212
$ws.Range("J13").Value = 0.640005219654821
$ws.Range("M13").Value = 12.167261
$ws.Range("N13").Value = 36.501783
$ws.Range("O13").Value = 0.0902104621952179
$ws.Range("P13").Value = 0.0902104621952179
$ws.Range("Q13").Value = 120.9220319476784
$ws.Range("R13").Value = 1088.298287529105
$ws.Range("S13").Value = 0.05773516667241337
$ws.Range("T13").Value = 0.05773516667241337

# Row 14
$ws.Range("G14").Value = 1.141526666666667
$ws.Range("H14").Value = 3.42458
$ws.Range("I14").Value = 0.0735117844504946
$ws.Range("J14").Value = 0.07351178445049458
$ws.Range("M14").Value = 43.97948166666666
$ws.Range("N14").Value = 131.938445
$ws.Range("O14").Value = 0.3260725128076164
$ws.Range("P14").Value = 0.3260725128076164
$ws.Range("Q14").Value = 50.20375110867777
$ws.Range("R14").Value = 451.8337599781
$ws.Range("S14").Value = 0.02397017227674464
$ws.Range("T14").Value = 0.02397017227674463

# Row 15
$ws.Range("G15").Value = 1.141526666666667
$ws.Range("H15").Value = 3.42458
$ws.Range("I15").Value = 0.0735117844504946
$ws.Range("J15").Value = 0.07351178445049458
$ws.Range("M15").Value = 57.80064033333333
$ws.Range("O15").Value = 0.4285452970598356
$ws.Range("P15").Value = 0.4285452970598356
$ws.Range("Q15").Value = 65.98097229090888
$ws.Range("R15").Value = 593.8287506181799
$ws.Range("S15").Value = 0.03150312950473581
$ws.Range("T15").Value = 0.03150312950473581

# Row 16
$ws.Range("G16").Value = 1.141526666666667
$ws.Range("H16").Value = 3.42458
$ws.Range("I16").Value = 0.0735117844504946
$ws.Range("J16").Value = 0.07351178445049458
$ws.Range("M16").Value = 20.92900166666667
$ws.Range("N16").Value = 62.787005
$ws.Range("O16").Value = 0.15517172793733
$ws.Range("P16").Value = 0.15517172793733
$ws.Range("Q16").Value = 23.89101350921111
$ws.Range("R16").Value = 215.0191215829
$ws.Range("S16").Value = 0.01140695061693979
$ws.Range("T16").Value = 0.01140695061693979

# Row 17
$ws.Range("G17").Value = 1.141526666666667
$ws.Range("H17").Value = 3.42458
$ws.Range("I17").Value = 0.0735117844504946
$ws.Range("J17").Value = 0.07351178445049458
$ws.Range("M17").Value = 12.167261
$ws.Range("N17").Value = 36.501783
$ws.Range("O17").Value = 0.0902104621952179
$ws.Range("P17").Value = 0.0902104621952179
$ws.Range("Q17").Value = 13.88925289179333
$ws.Range("R17").Value = 125.00327602614
$ws.Range("S17").Value = 0.00663153205207435
$ws.Range("T17").Value = 0.006631532052074349

Write-Output "Applied NATMI TPM update"